# Auto-generated edit script for LOM3218.docx restructuring
# Strategy: two-phase text move using unique placeholder tokens to avoid
# collisions between source/destination text during the big content reshuffle.
$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

$failures = New-Object System.Collections.ArrayList

# ---------- Phase 1: original text -> unique placeholder tokens ----------
$old0 = 'Apresentar aos alunos ingressantes o entendimento do que seja a carreira profissional e as bases conceituais da Engenharia Física, assim como empreendimentos na área.'
$ok0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_00@@", $wdReplaceAll)
if (-not $ok0) { [void]$failures.Add("phase1 idx=0 P6R1") }

$old1 = 'To introduce new students to an understanding of what a career is and the conceptual bases of Physical Engineering, as well as ventures in the area.'
$ok1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_01@@", $wdReplaceAll)
if (-not $ok1) { [void]$failures.Add("phase1 idx=1 P7R1") }

$old2 = '5817692 - Katia Cristiane Gandolpho Candioto' + [char]11 + ''
$ok2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_02@@", $wdReplaceAll)
if (-not $ok2) { [void]$failures.Add("phase1 idx=2 P9R1") }

$old3 = '1176388 - Luiz Tadeu Fernandes Eleno'
$ok3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_03@@", $wdReplaceAll)
if (-not $ok3) { [void]$failures.Add("phase1 idx=3 P9R2") }

$old4 = 'A carreira de Engenharia Física. Conceitos básicos de Engenharia. Competências e habilidades de um engenheiro. Física conceitual. Realização de experimentos e projetos de Engenharia Física.'
$ok4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_04@@", $wdReplaceAll)
if (-not $ok4) { [void]$failures.Add("phase1 idx=4 P11R1") }

$old5 = 'The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering.'
$ok5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_05@@", $wdReplaceAll)
if (-not $ok5) { [void]$failures.Add("phase1 idx=5 P12R1") }

$old6 = 'A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. ' + [char]11 + 'A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física.' + [char]11 + 'Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro.' + [char]11 + 'Desenvolvimento de um projeto temático de Engenharia Física.' + [char]11 + 'Competição entre projetos de diferentes grupos.' + [char]11 + 'Avaliação das competições e da disciplina como um todo.'
$ok6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_06@@", $wdReplaceAll)
if (-not $ok6) { [void]$failures.Add("phase1 idx=6 P14R1") }

$old7 = 'As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica).' + [char]11 + ''
$ok7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_07@@", $wdReplaceAll)
if (-not $ok7) { [void]$failures.Add("phase1 idx=7 P17R2") }

$old8 = 'A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo.' + [char]11 + ''
$ok8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_08@@", $wdReplaceAll)
if (-not $ok8) { [void]$failures.Add("phase1 idx=8 P17R4") }

$old9 = 'Devido às características da disciplina, não será oferecida recuperação.'
$ok9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_09@@", $wdReplaceAll)
if (-not $ok9) { [void]$failures.Add("phase1 idx=9 P17R6") }

$old10 = 'ARAÚJO-MOREIRA, F. M. Engenharia Física: a Carreira do Novo Milênio, São Carlos: Gráfica e Editora Guillen & Andriolli, 2014.' + [char]11 + 'BAZZO, A. B.; PEREIRA, L.T.V. Introdução à Engenharia. Editora da UFSC, Florianópolis, 1993.' + [char]11 + 'ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.' + [char]11 + 'BROCKMAN, J. B. Introdução à Engenharia. LTC, Rio de Janeiro, 2009.' + [char]11 + 'KNOWLEDGE FLOW. Engineering Physics - Ebook, Índia, 2015.' + [char]11 + 'CHAVES, A. S.; VALADARES, E. C.; ALVES, E. G. Aplicações da Física Quântica do Transistor à Nanotecnologia, São Paulo: Livraria da Física, 2005.'
$ok10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "@@MOVE_TOKEN_10@@", $wdReplaceAll)
if (-not $ok10) { [void]$failures.Add("phase1 idx=10 P19R1") }

# ---------- Phase 2: placeholder tokens -> final text ----------
$new0 = 'A carreira de Engenharia Física. Conceitos básicos de Engenharia. Competências e habilidades de um engenheiro. Física conceitual. Realização de experimentos e projetos de Engenharia Física.'
$ok2_0 = $d.Content.Find.Execute("@@MOVE_TOKEN_00@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new0, $wdReplaceAll)
if (-not $ok2_0) { [void]$failures.Add("phase2 idx=0 P6R1") }

$new1 = 'The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering.'
$ok2_1 = $d.Content.Find.Execute("@@MOVE_TOKEN_01@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new1, $wdReplaceAll)
if (-not $ok2_1) { [void]$failures.Add("phase2 idx=1 P7R1") }

$new2 = 'Apresentar aos alunos ingressantes o entendimento do que seja a carreira profissional e as bases conceituais da Engenharia Física, assim como empreendimentos na área.' + [char]11 + ''
$ok2_2 = $d.Content.Find.Execute("@@MOVE_TOKEN_02@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new2, $wdReplaceAll)
if (-not $ok2_2) { [void]$failures.Add("phase2 idx=2 P9R1") }

$new3 = 'A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. ' + [char]11 + 'A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física.' + [char]11 + 'Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro.' + [char]11 + 'Desenvolvimento de um projeto temático de Engenharia Física.' + [char]11 + 'Competição entre projetos de diferentes grupos.' + [char]11 + 'Avaliação das competições e da disciplina como um todo.'
$ok2_3 = $d.Content.Find.Execute("@@MOVE_TOKEN_03@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new3, $wdReplaceAll)
if (-not $ok2_3) { [void]$failures.Add("phase2 idx=3 P9R2") }

$new4 = 'As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica).'
$ok2_4 = $d.Content.Find.Execute("@@MOVE_TOKEN_04@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new4, $wdReplaceAll)
if (-not $ok2_4) { [void]$failures.Add("phase2 idx=4 P11R1") }

$new5 = 'To introduce new students to an understanding of what a career is and the conceptual bases of Physical Engineering, as well as ventures in the area.'
$ok2_5 = $d.Content.Find.Execute("@@MOVE_TOKEN_05@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new5, $wdReplaceAll)
if (-not $ok2_5) { [void]$failures.Add("phase2 idx=5 P12R1") }

$new6 = 'A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo.'
$ok2_6 = $d.Content.Find.Execute("@@MOVE_TOKEN_06@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new6, $wdReplaceAll)
if (-not $ok2_6) { [void]$failures.Add("phase2 idx=6 P14R1") }

$new7 = 'Devido às características da disciplina, não será oferecida recuperação.' + [char]11 + ''
$ok2_7 = $d.Content.Find.Execute("@@MOVE_TOKEN_07@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new7, $wdReplaceAll)
if (-not $ok2_7) { [void]$failures.Add("phase2 idx=7 P17R2") }

$new8 = 'ARAÚJO-MOREIRA, F. M. Engenharia Física: a Carreira do Novo Milênio, São Carlos: Gráfica e Editora Guillen & Andriolli, 2014.' + [char]11 + 'BAZZO, A. B.; PEREIRA, L.T.V. Introdução à Engenharia. Editora da UFSC, Florianópolis, 1993.' + [char]11 + 'ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.' + [char]11 + 'BROCKMAN, J. B. Introdução à Engenharia. LTC, Rio de Janeiro, 2009.' + [char]11 + 'KNOWLEDGE FLOW. Engineering Physics - Ebook, Índia, 2015.' + [char]11 + 'CHAVES, A. S.; VALADARES, E. C.; ALVES, E. G. Aplicações da Física Quântica do Transistor à Nanotecnologia, São Paulo: Livraria da Física, 2005.' + [char]11 + ''
$ok2_8 = $d.Content.Find.Execute("@@MOVE_TOKEN_08@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new8, $wdReplaceAll)
if (-not $ok2_8) { [void]$failures.Add("phase2 idx=8 P17R4") }

$new9 = '5817692 - Katia Cristiane Gandolpho Candioto'
$ok2_9 = $d.Content.Find.Execute("@@MOVE_TOKEN_09@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new9, $wdReplaceAll)
if (-not $ok2_9) { [void]$failures.Add("phase2 idx=9 P17R6") }

$new10 = '1176388 - Luiz Tadeu Fernandes Eleno'
$ok2_10 = $d.Content.Find.Execute("@@MOVE_TOKEN_10@@", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new10, $wdReplaceAll)
if (-not $ok2_10) { [void]$failures.Add("phase2 idx=10 P19R1") }

if ($failures.Count -gt 0) {
    Write-Output ("FAILURES: " + ($failures -join "; "))
} else {
    Write-Output "ALL OK"
}
